# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to reflect the newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates (row -> new value)
$wsExhibit.Range("F2").Value  = 6854
$wsExhibit.Range("F4").Value  = 441
$wsExhibit.Range("F12").Value = 38
$wsExhibit.Range("F14").Value = 431
$wsExhibit.Range("F16").Value = 1787
$wsExhibit.Range("F18").Value = 3473
$wsExhibit.Range("F21").Value = 18
$wsExhibit.Range("F22").Value = 2106
$wsExhibit.Range("F29").Value = 141

# 全部类型 sheet updates (row -> new value)
$wsAll.Range("F2").Value  = 6854
$wsAll.Range("F4").Value  = 441
$wsAll.Range("F13").Value = 38
$wsAll.Range("F15").Value = 431
$wsAll.Range("F17").Value = 1787
$wsAll.Range("F19").Value = 3473
$wsAll.Range("F22").Value = 18
$wsAll.Range("F23").Value = 2106
$wsAll.Range("F30").Value = 141
